# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal").
#
# A new week's data (3 rows: Especial/Primera/Segunda for Hayward kiwi,
# Región de O'Higgins) is inserted at the top of the dated price history
# (rows 20-22), pushing all subsequent rows down by 3 positions. The two
# oldest rows that fall off the bottom of the original range are replaced
# by two brand-new rows appended at the very end, keeping the same
# 3-rows-per-week cadence for that trailing week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture the existing data block (rows 20-127, columns A-T) before
#    shifting anything, then write it back starting 3 rows lower
#    (rows 23-130). This pushes the whole history down by 3 rows.
$srcBlock = $ws.Range("A20:T127")
$blockValues = $srcBlock.Value()
$destBlock = $ws.Range("A23:T130")
$destBlock.Value = $blockValues

# 2) Overwrite the newly-freed rows 20-22 with the new week's entries.
#    Only the date and the price columns (N, O, P, S) actually change -
#    everything else (market/region/product/variety/quality/volume/unit)
#    is identical to what was already there, so it's left untouched.
$newDates = New-Object 'object[,]' 3,1
$newDates[0,0] = 44687
$newDates[1,0] = 44687
$newDates[2,0] = 44687
$ws.Range("D20:D22").Value = $newDates

$newMinPrice = New-Object 'object[,]' 3,1
$newMinPrice[0,0] = 12000
$newMinPrice[1,0] = 10000
$newMinPrice[2,0] = 8000
$ws.Range("N20:N22").Value = $newMinPrice

$newMaxPrice = New-Object 'object[,]' 3,1
$newMaxPrice[0,0] = 12000
$newMaxPrice[1,0] = 10000
$newMaxPrice[2,0] = 8000
$ws.Range("O20:O22").Value = $newMaxPrice

$newAvgPrice = New-Object 'object[,]' 3,1
$newAvgPrice[0,0] = 12000
$newAvgPrice[1,0] = 10000
$newAvgPrice[2,0] = 8000
$ws.Range("P20:P22").Value = $newAvgPrice

$newKgPrice = New-Object 'object[,]' 3,1
$newKgPrice[0,0] = 667
$newKgPrice[1,0] = 556
$newKgPrice[2,0] = 444
$ws.Range("S20:S22").Value = $newKgPrice

# 3) The two brand-new trailing rows (129-130) sit beyond the sheet's
#    original used range, so the "Fecha" cells there pick up a default
#    date format instead of the column's usual one. Match the format
#    used throughout the rest of column D.
$ws.Range("D129:D130").NumberFormat = $ws.Range("D128").NumberFormat()
